# Trade #181 closed at 2026-02-17 22:04:21 - unknown UNKNOWN +0.000%
#
# This script mirrors the live-trading-bot's periodic workbook refresh:
#   1) Summary sheet roll-up numbers move (capital, P&L, trade counts, win rate).
#   2) Strategy Status row for "MarketMaking" moves in lock-step with Summary.
#   3) The open MarketMaking trade (Trade # 209) is closed out - fills in
#      Exit Price / Status / P&L / Capital After / Exit Reason / Duration -
#      on both the "All Trades" ledger and the per-strategy "MarketMaking" tab.
#   4) Two brand-new OPEN trades get appended (Trade # 242 on
#      "volatility_scorer", Trade # 243 on "MarketMaking") - again mirrored
#      onto "All Trades" and each strategy's own tab.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1) Summary sheet
# ---------------------------------------------------------------------------
$summary = $wb.Worksheets.Item("Summary")
$summary.Range("B3").Value = 1399.91    # Current Capital
$summary.Range("B4").Value = -0.31      # Total P&L $
$summary.Range("B5").Value = -0.03      # Total P&L %
$summary.Range("B6").Value = 209        # Total Trades
$summary.Range("B8").Value = 91         # Losing Trades
$summary.Range("B9").Value = 38.28      # Win Rate %

# ---------------------------------------------------------------------------
# 2) Strategy Status sheet - MarketMaking row (row 5)
# ---------------------------------------------------------------------------
$status = $wb.Worksheets.Item("Strategy Status")
$status.Range("C5").Value = 99.91    # Capital
$status.Range("D5").Value = 176      # Trades
$status.Range("E5").Value = -0.42    # P&L $
$status.Range("F5").Value = -0.09    # P&L %
$status.Range("G5").Value = 37.5     # Win Rate %

# ---------------------------------------------------------------------------
# 3) Close out Trade # 209 (MarketMaking) - "All Trades" ledger, row 210
#    Column order on "All Trades": ... K=Capital After, L=Exit Reason,
#    M=Duration (min), N=Entry Slippage, O=Exit Slippage ...
# ---------------------------------------------------------------------------
$allTrades = $wb.Worksheets.Item("All Trades")
$allTrades.Cells.Item(210, 7).Value = 0.176349      # G210 Exit Price
$allTrades.Cells.Item(210, 8).Value = "CLOSED"      # H210 Status
$allTrades.Cells.Item(210, 9).Value = -54.7824      # I210 P&L %
$allTrades.Cells.Item(210, 10).Value = -0.21        # J210 P&L $
$allTrades.Cells.Item(210, 11).Value = 99.91        # K210 Capital After
$allTrades.Cells.Item(210, 12).Value = "early_exit" # L210 Exit Reason
$allTrades.Cells.Item(210, 13).Value = 0.19         # M210 Duration (min)

# Same trade, mirrored on the "MarketMaking" strategy tab, row 177.
#    Column order here: ... K=Capital After, L=Entry Slippage, M=Exit
#    Slippage, N=Confidence, O=Entry Reason, P=Exit Reason, Q=Duration (min)
$marketMaking = $wb.Worksheets.Item("MarketMaking")
$marketMaking.Cells.Item(177, 7).Value = 0.176349      # G177 Exit Price
$marketMaking.Cells.Item(177, 8).Value = "CLOSED"      # H177 Status
$marketMaking.Cells.Item(177, 9).Value = -54.7824      # I177 P&L %
$marketMaking.Cells.Item(177, 10).Value = -0.21        # J177 P&L $
$marketMaking.Cells.Item(177, 11).Value = 99.91        # K177 Capital After
$marketMaking.Cells.Item(177, 16).Value = "early_exit" # P177 Exit Reason
$marketMaking.Cells.Item(177, 17).Value = 0.19         # Q177 Duration (min)

# ---------------------------------------------------------------------------
# 4a) New trade # 242 - volatility_scorer, appended to "All Trades" row 243
# ---------------------------------------------------------------------------
$allTrades.Cells.Item(243, 1).Value = 242
$allTrades.Cells.Item(243, 2).Value = "'2026-02-17"
$allTrades.Cells.Item(243, 3).Value = "22:04:13"
$allTrades.Cells.Item(243, 4).Value = "volatility_scorer"
$allTrades.Cells.Item(243, 5).Value = "NEUTRAL"
$allTrades.Cells.Item(243, 6).Value = 0.39
$allTrades.Cells.Item(243, 8).Value = "OPEN"
$allTrades.Cells.Item(243, 9).Value = 0
$allTrades.Cells.Item(243, 10).Value = 0
$allTrades.Cells.Item(243, 11).Value = 100
$allTrades.Cells.Item(243, 13).Value = 0
$allTrades.Cells.Item(243, 14).Value = 0
$allTrades.Cells.Item(243, 15).Value = 0
$allTrades.Cells.Item(243, 16).Value = 0.85
$allTrades.Cells.Item(243, 17).Value = "Low vol market (score: inf) - ideal for market making"

# Mirrored onto the "volatility_scorer" strategy tab, row 3.
$volScorer = $wb.Worksheets.Item("volatility_scorer")
$volScorer.Cells.Item(3, 1).Value = 242
$volScorer.Cells.Item(3, 2).Value = "'2026-02-17"
$volScorer.Cells.Item(3, 3).Value = "22:04:13"
$volScorer.Cells.Item(3, 4).Value = "volatility_scorer"
$volScorer.Cells.Item(3, 5).Value = "NEUTRAL"
$volScorer.Cells.Item(3, 6).Value = 0.39
$volScorer.Cells.Item(3, 8).Value = "OPEN"
$volScorer.Cells.Item(3, 9).Value = 0
$volScorer.Cells.Item(3, 10).Value = 0
$volScorer.Cells.Item(3, 11).Value = 100
$volScorer.Cells.Item(3, 12).Value = 0
$volScorer.Cells.Item(3, 13).Value = 0
$volScorer.Cells.Item(3, 14).Value = 0.85
$volScorer.Cells.Item(3, 15).Value = "Low vol market (score: inf) - ideal for market making"
$volScorer.Cells.Item(3, 17).Value = 0

# ---------------------------------------------------------------------------
# 4b) New trade # 243 - MarketMaking, appended to "All Trades" row 244
# ---------------------------------------------------------------------------
$allTrades.Cells.Item(244, 1).Value = 243
$allTrades.Cells.Item(244, 2).Value = "'2026-02-17"
$allTrades.Cells.Item(244, 3).Value = "22:04:14"
$allTrades.Cells.Item(244, 4).Value = "MarketMaking"
$allTrades.Cells.Item(244, 5).Value = "UP"
$allTrades.Cells.Item(244, 6).Value = 0.66075
$allTrades.Cells.Item(244, 8).Value = "OPEN"
$allTrades.Cells.Item(244, 9).Value = 0
$allTrades.Cells.Item(244, 10).Value = 0
$allTrades.Cells.Item(244, 11).Value = 100.1245541900307
$allTrades.Cells.Item(244, 13).Value = 0
$allTrades.Cells.Item(244, 14).Value = 0
$allTrades.Cells.Item(244, 15).Value = 0
$allTrades.Cells.Item(244, 16).Value = 0.6
$allTrades.Cells.Item(244, 17).Value = "Normal spread capture: 19600 bps"

# Mirrored onto the "MarketMaking" strategy tab, row 209.
$marketMaking.Cells.Item(209, 1).Value = 243
$marketMaking.Cells.Item(209, 2).Value = "'2026-02-17"
$marketMaking.Cells.Item(209, 3).Value = "22:04:14"
$marketMaking.Cells.Item(209, 4).Value = "MarketMaking"
$marketMaking.Cells.Item(209, 5).Value = "UP"
$marketMaking.Cells.Item(209, 6).Value = 0.66075
$marketMaking.Cells.Item(209, 8).Value = "OPEN"
$marketMaking.Cells.Item(209, 9).Value = 0
$marketMaking.Cells.Item(209, 10).Value = 0
$marketMaking.Cells.Item(209, 11).Value = 100.1245541900307
$marketMaking.Cells.Item(209, 12).Value = 0
$marketMaking.Cells.Item(209, 13).Value = 0
$marketMaking.Cells.Item(209, 14).Value = 0.6
$marketMaking.Cells.Item(209, 15).Value = "Normal spread capture: 19600 bps"
$marketMaking.Cells.Item(209, 17).Value = 0
